$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (date updated from 11-06 to 11-07)
$ws.Name = "Through 2021-11-07"

# Update the label in A12 (shared string text)
$ws.Range("A12").Value = "November (through 11-07)"

# Update November row (row 12) values
$ws.Range("B12").Value = 8
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 24
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 45
$ws.Range("H12").Value = 48

# Update Total row (row 13) values
$ws.Range("B13").Value = 266
$ws.Range("C13").Value = 503
$ws.Range("D13").Value = 734
$ws.Range("F13").Value = 492
$ws.Range("G13").Value = 1102
$ws.Range("H13").Value = 1492
